$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.121.75'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.46%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.107.58'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.17%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.88%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '350.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.21%  '

$ws.Range("E6").Value = '  -0.80%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5170'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.61%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4495'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '52.93'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.66%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08969'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.37%  '

$ws.Range("E11").Value = '  +0.65%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '25.81'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.24%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.104.56'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.69%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.773'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.21%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.176'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.69%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '99.58'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.26%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001152'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.10%  '

$ws.Range("E18").Value = '  -0.83%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '21.06'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +9.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06668'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.29%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.78%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.248'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.18%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.202.82'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.37%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.93'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.90%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.355'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.18%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.358.94'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.27%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.09'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.563'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.58%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '163.09'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.43%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.90'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.21%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.190'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.64%  '

$ws.Range("E32").Value = '  -0.34%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.654'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.71%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.280'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.23%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.963'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.14%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.947'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.34%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.23'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.32%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02595'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.50%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06859'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.51%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2319'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.02%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.57'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.23%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6843'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.15%  '

$ws.Range("E43").Value = '  +0.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.38'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.21%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6446'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.296'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.59%  '

$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000363'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.10%  '

$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.673'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.15%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '84.19'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.20%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.226'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.00%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07239'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.77%  '
